$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E. This shifts the existing
# input_position / input_projectTimeline / textarea_... columns one
# column to the right, preserving their values, styles and widths.
$ws.Columns("E").Insert()

# New column E header: "input_phone", formatted like the other header
# cells (bold, centered, thin border) so it reuses the sheet's existing
# bold/centered "Pandas" header style rather than creating a new one.
$ws.Range("E1").Value = "input_phone"
$ws.Range("E1").Font.Bold = $true
$ws.Range("E1").HorizontalAlignment = -4108
$ws.Range("E1").Borders.LineStyle = 1

# The sheet's column widths follow "header length + 2" characters
# (13 = len("input_phone") + 2). ColumnWidth is in "characters"; the
# engine adds a fixed 5/6-character pad when exporting to the stored
# <col width> attribute, so subtract it here to land exactly on 13.
$ws.Columns("E").ColumnWidth = 13 - 5/6

# Row 2 keeps a blank value under the new header, matching the other
# (pre-existing, empty) row-2 cells.
$ws.Range("E2").Value = ""
